$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert B32 into a text label ("dasf") and B33 into a static number (3),
# replacing the previous AVERAGE/STDEV formulas (labels are now strings).
$ws.Range("B32").Value = "dasf"
$ws.Range("B33").Value = 3

# Chart tweaks: show data labels to the right on both series and mark the
# value axis number format as no longer linked to the source cells.
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$ser1 = $chart.SeriesCollection().Item(1)
$ser2 = $chart.SeriesCollection().Item(2)
$ser1.DataLabels().Position = -4152
$ser2.DataLabels().Position = -4152
$valAxis = $chart.Axes(2)
$valAxis.TickLabels.NumberFormatLinked = 0

# Keep the saved selection consistent with the authored workbook.
$null = $ws.Range("B33").Select()
